# Try Convert Feature Login
# - edit name keyword > keywords
# - add folder for test web
# - add keyword for teardown login
# - edit name capture screen shot
# - add "sleep 2" to keyword Capture screen shot
# - add localized
# - add repository
# - add test case (converted) 1-9
#
# Applies the test-case content updates to the "Login" QA traceback sheet:
#   * row 4 / row 5 now carry the same login credentials note as row 3
#   * the 7th test case (row 7) is reclassified Regression -> Non-Regression
#   * several rows grew taller to fit their wrapped text
#   * the view scrolled down to show the newly edited rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

# --- Row 4 & Row 5: fill in the missing "Source-Requisite" test data ----
# (same credentials text already used in row 3 / column C)
$ws.Range("C4").Value = "Username:admin`nPassword:12345678"
$ws.Range("C4").WrapText = $true

$ws.Range("C5").Value = "Username:admin`nPassword:12345678"
$ws.Range("C5").WrapText = $true

# --- Row 7: reclassify the test case from Regression to Non-Regression --
$ws.Range("G7").Value = "Non-Regression"

# --- Row heights: rows grew to fit wrapped content -----------------------
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 90
$ws.Rows.Item(9).RowHeight = 90
$ws.Rows.Item(10).RowHeight = 90
$ws.Rows.Item(11).RowHeight = 90

# --- Sheet view moved down, selection now on E11 --------------------------
$ws.Range("E11").Select()

# --- Print setup: force portrait orientation -------------------------------
$ws.PageSetup.Orientation = 1
